# Generate Report for handoff
# - Removes the "a49e635b-..." row (handed-back/in-sync row) from every sheet,
#   shifting the ".localization-config" row up by one.
# - Updates the status text for the remaining "10fb5430-..." row from
#   "Handed back: in sync with en-US" to "Ready for handoff".
# - Bumps the "Latest Handoff Datetime" for the remaining row (new handoff run).
# - Rebuilds the hyperlinks collection on every sheet to match the shifted rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Drop row 3 (a49e635b...), row 4 (.localization-config) shifts up to row 3.
$ws1.Rows(3).EntireRow.Delete()

# Update status text for the 10fb5430 row.
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

# Rebuild hyperlinks (old refs became stale after the row delete).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7b5d1cc19b982b6c4867d75646419aa94e477968/e2e/10fb5430-6c98-4264-9ae7-b22d2f391eab.md", [Type]::Missing, [Type]::Missing, "10fb5430-6c98-4264-9ae7-b22d2f391eab.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7b5d1cc19b982b6c4867d75646419aa94e477968/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# Restore the name-column styling (hyperlink-add resets the cell style).
$ws1.Range("A2").Style = "HyperLink"
$ws1.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows(3).EntireRow.Delete()

$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "2016-01-25 13:10:26"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7b5d1cc19b982b6c4867d75646419aa94e477968/e2e/10fb5430-6c98-4264-9ae7-b22d2f391eab.md", [Type]::Missing, [Type]::Missing, "10fb5430-6c98-4264-9ae7-b22d2f391eab.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f39d4cf2d2ffd6c97b188433fb8426c6ac063f0d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/10fb5430-6c98-4264-9ae7-b22d2f391eab.4b74a62e6652c4a47efe9b32d19fe2a955ed42a4.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "10fb5430-6c98-4264-9ae7-b22d2f391eab.4b74a62e6652c4a47efe9b32d19fe2a955ed42a4.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b1037fa410a88f258a8a60904175bffadf443f60/e2e/10fb5430-6c98-4264-9ae7-b22d2f391eab.md", [Type]::Missing, [Type]::Missing, "10fb5430-6c98-4264-9ae7-b22d2f391eab.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/70cc8944ade2b49f392ddb50b9d9b98c6c090749/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/10fb5430-6c98-4264-9ae7-b22d2f391eab.4b74a62e6652c4a47efe9b32d19fe2a955ed42a4.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "10fb5430-6c98-4264-9ae7-b22d2f391eab.4b74a62e6652c4a47efe9b32d19fe2a955ed42a4.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7b5d1cc19b982b6c4867d75646419aa94e477968/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

$ws2.Range("A2").Style = "HyperLink"
$ws2.Range("C2").Style = "HyperLink"
$ws2.Range("E2").Style = "HyperLink"
$ws2.Range("F2").Style = "HyperLink"
$ws2.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows(3).EntireRow.Delete()

$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "2016-01-25 13:10:36"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7b5d1cc19b982b6c4867d75646419aa94e477968/e2e/10fb5430-6c98-4264-9ae7-b22d2f391eab.md", [Type]::Missing, [Type]::Missing, "10fb5430-6c98-4264-9ae7-b22d2f391eab.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f2c696f9bf141bb08014e48b79eeb149015b62a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/10fb5430-6c98-4264-9ae7-b22d2f391eab.4b74a62e6652c4a47efe9b32d19fe2a955ed42a4.de-de.xlf", [Type]::Missing, [Type]::Missing, "10fb5430-6c98-4264-9ae7-b22d2f391eab.4b74a62e6652c4a47efe9b32d19fe2a955ed42a4.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1f19895180e377555938d9d5f27232b82cb1f2d1/e2e/10fb5430-6c98-4264-9ae7-b22d2f391eab.md", [Type]::Missing, [Type]::Missing, "10fb5430-6c98-4264-9ae7-b22d2f391eab.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c762b6590987ad924028f2da652b5a17306deef5/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/10fb5430-6c98-4264-9ae7-b22d2f391eab.4b74a62e6652c4a47efe9b32d19fe2a955ed42a4.de-de.xlf", [Type]::Missing, [Type]::Missing, "10fb5430-6c98-4264-9ae7-b22d2f391eab.4b74a62e6652c4a47efe9b32d19fe2a955ed42a4.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7b5d1cc19b982b6c4867d75646419aa94e477968/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

$ws3.Range("A2").Style = "HyperLink"
$ws3.Range("C2").Style = "HyperLink"
$ws3.Range("E2").Style = "HyperLink"
$ws3.Range("F2").Style = "HyperLink"
$ws3.Range("A3").Style = "HyperLink"
